$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B/C column header labels (CA_FFR <-> CA_LF)
$ws.Range("B1").Value = "CA_LF"
$ws.Range("C1").Value = "CA_FFR"

# Swap the B2/C2 numeric values
$ws.Range("B2").Value = 0.8376754650626033
$ws.Range("C2").Value = 10.56400606921197
